# Experiment 222 / Tabelle_3: correct the heating-current (C) and
# heating-voltage-error (D) formulas for rows 2-13.
#
# The original formulas multiplied by an extra *5 (C: "=X*5",
# D: "=C{row}*0.008+0.03*5"). That factor was a mistake and is removed:
#   C{row}: "=X"                    (was "=X*5")
#   D{row}: "=C{row}*0.008+0.03"    (was "=C{row}*0.008+0.03*5")
#
# Column C literal multipliers per row (from the original formulas):
#   2:2.62  3:2.61  4:2.62  5:2.62  6:2.62  7:2.62
#   8:2.61  9:2.62 10:2.61 11:2.61 12:2.61 13:2.61

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cValues = @{
    2  = "2.62"
    3  = "2.61"
    4  = "2.62"
    5  = "2.62"
    6  = "2.62"
    7  = "2.62"
    8  = "2.61"
    9  = "2.62"
    10 = "2.61"
    11 = "2.61"
    12 = "2.61"
    13 = "2.61"
}

for ($row = 2; $row -le 13; $row++) {
    $ws.Range("C$row").Formula = "=" + $cValues[$row]
    $ws.Range("D$row").Formula = "=C$row*0.008+0.03"
}

# The saved file shows D2 as the active selection.
$ws.Range("D2").Select()
